$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2531645569620253
$ws.Range("C2").Value = 0.4556962025316456
$ws.Range("J2").Value = 0.02531645569620253
$ws.Range("P2").Value = 0.189873417721519
$ws.Range("S2").Value = 0.0759493670886076
$ws.Range("C3").Value = 0.02702702702702703
$ws.Range("J3").Value = 0.05405405405405406
$ws.Range("P3").Value = 0.8378378378378378
$ws.Range("S3").Value = 0.08108108108108109
$ws.Range("J4").Value = 0
$ws.Range("P4").Value = 0.6363636363636364
$ws.Range("S4").Value = 0.3636363636363636
$ws.Range("B6").Value = 0.02439024390243903
$ws.Range("J6").Value = 0.1951219512195122
$ws.Range("O6").Value = 0.02439024390243903
$ws.Range("Q6").Value = 0.3170731707317073
$ws.Range("R6").Value = 0.1219512195121951
$ws.Range("S6").Value = 0.3170731707317073
$ws.Range("B7").Value = 0.03125
$ws.Range("D7").Value = 0.03125
$ws.Range("F7").Value = 0.03125
$ws.Range("J7").Value = 0.09375
$ws.Range("Q7").Value = 0.15625
$ws.Range("R7").Value = 0.15625
$ws.Range("S7").Value = 0.5
$ws.Range("B8").Value = 0.06122448979591837
$ws.Range("D8").Value = 0.03061224489795918
$ws.Range("F8").Value = 0.07142857142857142
$ws.Range("J8").Value = 0.2244897959183673
$ws.Range("O8").Value = 0.03061224489795918
$ws.Range("Q8").Value = 0.2448979591836735
$ws.Range("R8").Value = 0.1122448979591837
$ws.Range("S8").Value = 0.2244897959183673
$ws.Range("B9").Value = 0.06521739130434782
$ws.Range("D9").Value = 0.04347826086956522
$ws.Range("F9").Value = 0.04347826086956522
$ws.Range("J9").Value = 0.1739130434782609
$ws.Range("Q9").Value = 0.2826086956521739
$ws.Range("R9").Value = 0.1521739130434783
$ws.Range("S9").Value = 0.2391304347826087
$ws.Range("B10").Value = 0.09533898305084745
$ws.Range("D10").Value = 0.01059322033898305
$ws.Range("E10").Value = 0.00211864406779661
$ws.Range("F10").Value = 0.0423728813559322
$ws.Range("J10").Value = 0.1334745762711864
$ws.Range("O10").Value = 0.00423728813559322
$ws.Range("Q10").Value = 0.3453389830508475
$ws.Range("R10").Value = 0.1101694915254237
$ws.Range("S10").Value = 0.2563559322033898
$ws.Range("G11").Value = 0.1666666666666667
$ws.Range("J11").Value = 0.1
$ws.Range("K11").Value = 0.1833333333333333
$ws.Range("L11").Value = 0.5333333333333333
$ws.Range("S11").Value = 0.01666666666666667
$ws.Range("G12").Value = 0.625
$ws.Range("J12").Value = 0.3125
$ws.Range("L12").Value = 0.03125
$ws.Range("S12").Value = 0.03125
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.5
$ws.Range("H15").Value = 0.09615384615384616
$ws.Range("I15").Value = 0.03846153846153846
$ws.Range("J15").Value = 0.5
$ws.Range("K15").Value = 0.0576923076923077
$ws.Range("M15").Value = 0.01923076923076923
$ws.Range("O15").Value = 0.09615384615384616
$ws.Range("S15").Value = 0.1923076923076923
$ws.Range("H16").Value = 0.07692307692307693
$ws.Range("I16").Value = 0.09615384615384616
$ws.Range("J16").Value = 0.5576923076923077
$ws.Range("K16").Value = 0.0576923076923077
$ws.Range("S16").Value = 0.2115384615384615
$ws.Range("F17").Value = 0.009389671361502348
$ws.Range("H17").Value = 0.1690140845070423
$ws.Range("I17").Value = 0.08450704225352113
$ws.Range("J17").Value = 0.5446009389671361
$ws.Range("K17").Value = 0.06572769953051644
$ws.Range("N17").Value = 0.004694835680751174
$ws.Range("O17").Value = 0.03286384976525822
$ws.Range("S17").Value = 0.0892018779342723
$ws.Range("F18").Value = 0.01265822784810127
$ws.Range("H18").Value = 0.1012658227848101
$ws.Range("I18").Value = 0.0759493670886076
$ws.Range("J18").Value = 0.5443037974683544
$ws.Range("K18").Value = 0.06329113924050633
$ws.Range("O18").Value = 0.0759493670886076
$ws.Range("S18").Value = 0.1265822784810127
$ws.Range("F19").Value = 0.007380073800738007
$ws.Range("H19").Value = 0.1660516605166052
$ws.Range("I19").Value = 0.05535055350553506
$ws.Range("J19").Value = 0.5092250922509225
$ws.Range("K19").Value = 0.08856088560885608
$ws.Range("M19").Value = 0.003690036900369004
$ws.Range("O19").Value = 0.07749077490774908
$ws.Range("S19").Value = 0.09225092250922509
